$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Canada predictions: append 3 new weekly rows (47-49) for prediction date
# 2021-01-09, mirroring the existing rows 44-46 pattern but adding the new
# "24 Jan -- 30 Jan 2021" target week.
#
# Column A holds an ISO-looking date string ("2021-01-09") that must stay
# TEXT (as it already is for every other row in the sheet) rather than be
# auto-converted to a date serial by Excel's smart entry. Copying the cell
# from an existing row that already holds that exact text preserves its
# text-ness without touching number formats / styles.
$ws.Range("A46").Copy($ws.Range("A47"))
$ws.Range("A46").Copy($ws.Range("A48"))
$ws.Range("A46").Copy($ws.Range("A49"))

$ws.Cells.Item(47, 2).Value = "10 Jan -- 16 Jan 2021"
$ws.Cells.Item(47, 4).Value = 124.04
$ws.Cells.Item(47, 6).Value = "KNN"

$ws.Cells.Item(48, 2).Value = "17 Jan -- 23 Jan 2021"
$ws.Cells.Item(48, 4).Value = 114.28
$ws.Cells.Item(48, 6).Value = "KNN"

$ws.Cells.Item(49, 2).Value = "24 Jan -- 30 Jan 2021"
$ws.Cells.Item(49, 4).Value = 120.01
$ws.Cells.Item(49, 6).Value = "KNN"
